$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 3 - Texas -- Bexar County
$ws.Range("B3").Value = 44027
$ws.Range("C3").Value = 21546
$ws.Range("D3").Value = 208

# Row 9 - Arkansas
$ws.Range("C9").Value = 30297
$ws.Range("D9").Value = 335
$ws.Range("E9").Value = 6444
$ws.Range("F9").Value = 87
$ws.Range("G9").Value = 24.63
$ws.Range("K9").Value = 26159
$ws.Range("L9").Value = 333

# Row 10 - California - San Diego
$ws.Range("B10").Value = 44027
$ws.Range("C10").Value = 21446
$ws.Range("D10").Value = 448
$ws.Range("E10").Value = 786
$ws.Range("F10").Value = 18
$ws.Range("G10").Value = 4.67
$ws.Range("H10").Value = 4.14
$ws.Range("K10").Value = 16838
$ws.Range("L10").Value = 435

# Row 25 - Nebraska
$ws.Range("B25").Value = 44027
$ws.Range("C25").Value = 21979
$ws.Range("D25").Value = 291
$ws.Range("E25").Value = 1286
$ws.Range("G25").Value = 7.61
$ws.Range("H25").Value = 8
$ws.Range("K25").Value = 16891
$ws.Range("L25").Value = 275

# Row 27 - California
$ws.Range("B27").NumberFormat = "YYYY-MM-DD"
$ws.Range("B27").Value = 44026
$ws.Range("C27").Value = 347634
$ws.Range("D27").Value = 7164
$ws.Range("E27").Value = 9697
$ws.Range("F27").Value = 621
$ws.Range("G27").Value = 4.3
$ws.Range("H27").Value = 8.800000000000001
$ws.Range("J27").Value = $true
$ws.Range("K27").Value = 223576
$ws.Range("L27").Value = 7038
$ws.Range("O27").Value = "Success!"

# Row 32 - Washington
$ws.Range("B32").Value = 44027
$ws.Range("C32").Value = 43046
$ws.Range("D32").Value = 1421
$ws.Range("E32").Value = 1649
$ws.Range("G32").Value = 5.45
$ws.Range("H32").Value = 3.53
$ws.Range("K32").Value = 30278
$ws.Range("L32").Value = 1359

# Row 36 - Iowa
$ws.Range("C36").Value = 36322
$ws.Range("D36").Value = 772
$ws.Range("E36").Value = 3104
$ws.Range("F36").Value = 37
$ws.Range("H36").Value = 4.79
